$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H28").Value = 749.375
$ws.Range("I28").Value = 408.18182
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 408.18182
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = 76.81817999999998
$ws.Range("N28").Value = -2470
$ws.Range("H33").Value = 1310.1875
$ws.Range("I33").Value = 213.22223
$ws.Range("J33").Value = 2720.5715
$ws.Range("K33").Value = 213.22223
$ws.Range("L33").Value = 2720.5715
$ws.Range("M33").Value = 15.77777
$ws.Range("N33").Value = -3178.5715
$ws.Range("H74").Value = 4996.1763
$ws.Range("I74").Value = 4343.6
$ws.Range("J74").Value = 5928.4287
$ws.Range("K74").Value = 4343.6
$ws.Range("L74").Value = 5928.4287
$ws.Range("M74").Value = -3407.6
$ws.Range("N74").Value = -7800.4287
$ws.Range("H77").Value = 4996.1763
$ws.Range("I77").Value = 4343.6
$ws.Range("J77").Value = 5928.4287
$ws.Range("K77").Value = 21718
$ws.Range("L77").Value = 29642.1435
$ws.Range("M77").Value = -17038
$ws.Range("N77").Value = -39002.14350000001
$ws.Range("H129").Value = 909.14545
$ws.Range("I129").Value = 684.8
$ws.Range("J129").Value = 959
$ws.Range("K129").Value = 2054.4
$ws.Range("L129").Value = 2877
$ws.Range("M129").Value = 2945.6
$ws.Range("N129").Value = -12877
$ws.Range("H137").Value = 1330.2363
$ws.Range("I137").Value = 827.7778
$ws.Range("J137").Value = 1574.6757
$ws.Range("K137").Value = 2483.3334
$ws.Range("L137").Value = 4724.0271
$ws.Range("M137").Value = 66.66660000000002
$ws.Range("N137").Value = -9824.027099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32889.465
$ws.Range("I32").Value = 6781.375
$ws.Range("J32").Value = 98159.69
$ws.Range("K32").Value = 6781.375
$ws.Range("L32").Value = 98159.69
$ws.Range("M32").Value = -6494.375
$ws.Range("H45").Value = 1796
$ws.Range("I45").Value = 700
$ws.Range("J45").Value = 2070
$ws.Range("K45").Value = 700
$ws.Range("L45").Value = 2070
$ws.Range("M45").Value = -323
$ws.Range("H74").Value = 1659.0857
$ws.Range("I74").Value = 1232.5555
$ws.Range("J74").Value = 2110.7058
$ws.Range("K74").Value = 1232.5555
$ws.Range("L74").Value = 2110.7058
$ws.Range("M74").Value = -358.5554999999999
$ws.Range("N74").Value = -3858.7058
$ws.Range("H77").Value = 1659.0857
$ws.Range("I77").Value = 1232.5555
$ws.Range("J77").Value = 2110.7058
$ws.Range("K77").Value = 6162.7775
$ws.Range("L77").Value = 10553.529
$ws.Range("M77").Value = -1794.7775
$ws.Range("N77").Value = -19289.529
$ws.Range("H102").Value = 52040.7
$ws.Range("I102").Value = 85202.914
$ws.Range("J102").Value = 2297.375
$ws.Range("K102").Value = 85202.914
$ws.Range("L102").Value = 2297.375
$ws.Range("M102").Value = -83580.914
$ws.Range("N102").Value = -5541.375
$ws.Range("H122").Value = 2118.238
$ws.Range("I122").Value = 1872.0714
$ws.Range("J122").Value = 2610.5715
$ws.Range("K122").Value = 5616.2142
$ws.Range("L122").Value = 7831.7145
$ws.Range("M122").Value = -3166.2142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 521.25
$ws.Range("I12").Value = 521.25
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 521.25
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -353.25
$ws.Range("H26").Value = 17388.2
$ws.Range("I26").Value = 8235.25
$ws.Range("J26").Value = 54000
$ws.Range("K26").Value = 8235.25
$ws.Range("L26").Value = 54000
$ws.Range("M26").Value = -7943.25
$ws.Range("N26").Value = -54584
$ws.Range("H28").Value = 29800
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 29800
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 29800
$ws.Range("N28").Value = -30388
$ws.Range("H42").Value = 398000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 398000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 398000
$ws.Range("N42").Value = -398656
$ws.Range("H105").Value = 92864
$ws.Range("I105").Value = 101779.6
$ws.Range("J105").Value = 85434.336
$ws.Range("K105").Value = 101779.6
$ws.Range("L105").Value = 85434.336
$ws.Range("M105").Value = -100032.6
$ws.Range("N105").Value = -88928.336
$ws.Range("H134").Value = 1769.36
$ws.Range("I134").Value = 1746.3636
$ws.Range("J134").Value = 1938
$ws.Range("K134").Value = 5239.0908
$ws.Range("L134").Value = 5814
$ws.Range("M134").Value = -2704.0908
$ws.Range("N134").Value = -10884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1861.2766
$ws.Range("I31").Value = 1230.7916
$ws.Range("J31").Value = 2519.1738
$ws.Range("K31").Value = 1230.7916
$ws.Range("L31").Value = 2519.1738
$ws.Range("M31").Value = -935.7916
$ws.Range("N31").Value = -3109.1738
$ws.Range("H34").Value = 1861.2766
$ws.Range("I34").Value = 1230.7916
$ws.Range("J34").Value = 2519.1738
$ws.Range("K34").Value = 1230.7916
$ws.Range("L34").Value = 2519.1738
$ws.Range("M34").Value = -1028.7916
$ws.Range("N34").Value = -2923.1738
$ws.Range("H62").Value = 6947420
$ws.Range("I62").Value = 18520686
$ws.Range("J62").Value = 3460
$ws.Range("K62").Value = 18520686
$ws.Range("L62").Value = 3460
$ws.Range("M62").Value = -18520062
$ws.Range("N62").Value = -4708
$ws.Range("H65").Value = 6947420
$ws.Range("I65").Value = 18520686
$ws.Range("J65").Value = 3460
$ws.Range("K65").Value = 92603430
$ws.Range("L65").Value = 17300
$ws.Range("M65").Value = -92600310
$ws.Range("N65").Value = -23540
$ws.Range("H107").Value = 1064.6364
$ws.Range("I107").Value = 2116.6667
$ws.Range("J107").Value = 670.125
$ws.Range("K107").Value = 2116.6667
$ws.Range("L107").Value = 670.125
$ws.Range("M107").Value = -196.6667000000002
$ws.Range("N107").Value = -4510.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 273146.62
$ws.Range("I107").Value = 472.10345
$ws.Range("J107").Value = 632581.25
$ws.Range("K107").Value = 1416.31035
$ws.Range("L107").Value = 1897743.75
$ws.Range("M107").Value = 503.68965
$ws.Range("N107").Value = -1901583.75
$ws.Range("H131").Value = 913.68335
$ws.Range("I131").Value = 515.4
$ws.Range("J131").Value = 949.8909
$ws.Range("K131").Value = 1546.2
$ws.Range("L131").Value = 2849.6727
$ws.Range("M131").Value = 3493.8
$ws.Range("N131").Value = -12929.6727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H102").Value = 2056.7273
$ws.Range("I102").Value = 2112.4
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 2112.4
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -490.4000000000001
$ws.Range("N102").Value = -4744
$ws.Range("H113").Value = 1407.75
$ws.Range("I113").Value = 985.1667
$ws.Range("J113").Value = 1830.3334
$ws.Range("K113").Value = 985.1667
$ws.Range("L113").Value = 1830.3334
$ws.Range("M113").Value = 1184.8333
$ws.Range("N113").Value = -6170.3334
$ws.Range("H126").Value = 4218.4
$ws.Range("I126").Value = 5364
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 16092
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -13622
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1506.375
$ws.Range("I7").Value = 1275.2858
$ws.Range("J7").Value = 3124
$ws.Range("K7").Value = 1275.2858
$ws.Range("L7").Value = 3124
$ws.Range("M7").Value = -1163.2858
$ws.Range("N7").Value = -3348
$ws.Range("H82").Value = 2659.889
$ws.Range("I82").Value = 1795.8
$ws.Range("J82").Value = 3740
$ws.Range("K82").Value = 1795.8
$ws.Range("L82").Value = 3740
$ws.Range("M82").Value = -1434.8
$ws.Range("H85").Value = 2659.889
$ws.Range("I85").Value = 1795.8
$ws.Range("J85").Value = 3740
$ws.Range("K85").Value = 1795.8
$ws.Range("L85").Value = 3740
$ws.Range("M85").Value = -547.8
$ws.Range("H96").Value = 13891.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 13891.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 13891.5
$ws.Range("N96").Value = -19383.5
$ws.Range("H100").Value = 2698.75
$ws.Range("I100").Value = 2155.1428
$ws.Range("J100").Value = 3459.8
$ws.Range("K100").Value = 2155.1428
$ws.Range("L100").Value = 3459.8
$ws.Range("M100").Value = -1614.1428
$ws.Range("N100").Value = -4541.8
$ws.Range("H126").Value = 1506.375
$ws.Range("I126").Value = 1275.2858
$ws.Range("J126").Value = 3124
$ws.Range("K126").Value = 3825.8574
$ws.Range("L126").Value = 9372
$ws.Range("M126").Value = -1355.8574
$ws.Range("N126").Value = -14312
$ws.Range("H132").Value = 3111.9
$ws.Range("I132").Value = 3373.0334
$ws.Range("J132").Value = 2720.2
$ws.Range("K132").Value = 10119.1002
$ws.Range("L132").Value = 8160.599999999999
$ws.Range("M132").Value = -7589.100199999999
$ws.Range("N132").Value = -13220.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 400659.8
$ws.Range("I81").Value = 333833.34
$ws.Range("J81").Value = 500899.5
$ws.Range("K81").Value = 667666.6800000001
$ws.Range("L81").Value = 1001799
$ws.Range("M81").Value = -666605.6800000001
$ws.Range("N81").Value = -1003921
$ws.Range("H84").Value = 400659.8
$ws.Range("I84").Value = 333833.34
$ws.Range("J84").Value = 500899.5
$ws.Range("K84").Value = 3338333.4
$ws.Range("L84").Value = 5008995
$ws.Range("M84").Value = -3333029.4
$ws.Range("N84").Value = -5019603
$ws.Range("H113").Value = 819.53845
$ws.Range("I113").Value = 449.75
$ws.Range("J113").Value = 983.8889
$ws.Range("K113").Value = 1349.25
$ws.Range("L113").Value = 2951.6667
$ws.Range("M113").Value = 820.75
$ws.Range("H119").Value = 44990
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 44990
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 44990
$ws.Range("N119").Value = -54666
$ws.Range("H122").Value = 1549.9166
$ws.Range("I122").Value = 1439.9
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 4319.700000000001
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -1869.700000000001
$ws.Range("N122").Value = -11200
$ws.Range("H126").Value = 1890.9333
$ws.Range("I126").Value = 1717.6666
$ws.Range("J126").Value = 2150.8333
$ws.Range("K126").Value = 5152.9998
$ws.Range("L126").Value = 6452.499899999999
$ws.Range("M126").Value = -2682.9998
$ws.Range("H132").Value = 1501.9487
$ws.Range("I132").Value = 1056.1154
$ws.Range("J132").Value = 2393.6155
$ws.Range("K132").Value = 3168.3462
$ws.Range("L132").Value = 7180.8465
$ws.Range("M132").Value = -638.3462
$ws.Range("N132").Value = -12240.8465
$ws.Range("H136").Value = 1411.625
$ws.Range("I136").Value = 735.1818
$ws.Range("J136").Value = 2899.8
$ws.Range("K136").Value = 2205.5454
$ws.Range("L136").Value = 8699.400000000001
$ws.Range("M136").Value = 344.4546
$ws.Range("N136").Value = -13799.4
